# saibabacharityreceiptor/template.xlsx
# "improved css and excel download"
#
# - K1/L1 header labels get explicit date-format hints appended.
# - Columns K and L are widened so the longer labels are readable.
# - The saved view no longer has a scrolled topLeftCell, and the
#   remembered selection moves from A2 to L8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
$ws.Range("K1").Value = "Date Received (mm/dd/yyyy)"
$ws.Range("L1").Value = "Issued Date(mm/dd/yyyy)"

# --- Column width updates --------------------------------------------------
# ColumnWidth is expressed in "characters"; Excel stores width on disk as
# characters + 5/6. Back that padding out so the saved width lands on the
# target values (25.5703125 and 21).
$ws.Columns.Item(11).ColumnWidth = 25.5703125 - 5/6
$ws.Columns.Item(12).ColumnWidth = 21 - 5/6

# --- Selection / scroll position -------------------------------------------
$ws.Range("L8").Select() | Out-Null
